$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.509218573570251
$ws.Range("B1").Value = 1.359439134597778
$ws.Range("C1").Value = 4.554666042327881
$ws.Range("D1").Value = 2.08967113494873
$ws.Range("E1").Value = 0.7075942754745483
